$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# CC_cal_temp
$ws.Range("H6").Value = 20.05

# CC_wl
$ws.Range("H7").Value = '[189.68, 190.48, 191.27, 192.07, 192.86, 193.66, 194.45, 195.25, 196.04, 196.84, 197.64, 198.43, 199.23, 200.03, 200.82, 201.62, 202.42, 203.22, 204.01, 204.81, 205.61, 206.41, 207.21, 208.01, 208.8, 209.6, 210.4, 211.2, 212, 212.8, 213.6, 214.4, 215.2, 216, 216.81, 217.61, 218.41, 219.21, 220.01, 220.81, 221.62, 222.42, 223.22, 224.02, 224.82, 225.63, 226.43, 227.23, 228.04, 228.84, 229.64, 230.45, 231.25, 232.05, 232.86, 233.66, 234.47, 235.27, 236.08, 236.88, 237.69, 238.49, 239.3, 240.1, 240.91, 241.71, 242.52, 243.32, 244.13, 244.94, 245.74, 246.55, 247.35, 248.16, 248.97, 249.77, 250.58, 251.39, 252.19, 253, 253.81, 254.62, 255.42, 256.23, 257.04, 257.85, 258.65, 259.46, 260.27, 261.08, 261.88, 262.69, 263.5, 264.31, 265.12, 265.92, 266.73, 267.54, 268.35, 269.16, 269.97, 270.77, 271.58, 272.39, 273.2, 274.01, 274.82, 275.63, 276.44, 277.25, 278.05, 278.86, 279.67, 280.48, 281.29, 282.1, 282.91, 283.72, 284.53, 285.34, 286.15, 286.95, 287.76, 288.57, 289.38, 290.19, 291, 291.81, 292.62, 293.43, 294.24, 295.05, 295.86, 296.67, 297.48, 298.28, 299.09, 299.9, 300.71, 301.52, 302.33, 303.14, 303.95, 304.76, 305.57, 306.38, 307.19, 307.99, 308.8, 309.61, 310.42, 311.23, 312.04, 312.85, 313.66, 314.47, 315.27, 316.08, 316.89, 317.7, 318.51, 319.32, 320.12, 320.93, 321.74, 322.55, 323.36, 324.17, 324.97, 325.78, 326.59, 327.4, 328.2, 329.01, 329.82, 330.63, 331.43, 332.24, 333.05, 333.85, 334.66, 335.47, 336.28, 337.08, 337.89, 338.69, 339.5, 340.31, 341.11, 341.92, 342.73, 343.53, 344.34, 345.14, 345.95, 346.75, 347.56, 348.36, 349.17, 349.97, 350.78, 351.58, 352.39, 353.19, 353.99, 354.8, 355.6, 356.41, 357.21, 358.01, 358.82, 359.62, 360.42, 361.23, 362.03, 362.83, 363.63, 364.43, 365.24, 366.04, 366.84, 367.64, 368.44, 369.24, 370.05, 370.85, 371.65, 372.45, 373.25, 374.05, 374.85, 375.65, 376.45, 377.25, 378.04, 378.84, 379.64, 380.44, 381.24, 382.04, 382.84, 383.63, 384.43, 385.23, 386.02, 386.82, 387.62, 388.41, 389.21, 390.01, 390.8, 391.6, 392.39, 393.19, 393.98, 394.78]'

# CC_eno3
$ws.Range("H8").Value = '[-0.00024427, -0.00024148, -0.00023872, -0.00221546, 0.01922061, -0.00388358, -0.00190129, -0.00538667, -0.00394554, -0.00227395, -0.00072157, 0.00469364, -0.00996021, 0.01878355, 0.00506595, -0.00020259, -0.0001998, 0.00060761, -0.00019425, -0.00019146, 0.00117028, 0.01003073, 0.00114032, 0.00560789, 0.00695074, 0.00767565, 0.00742252, 0.00735185, 0.00683618, 0.00657612, 0.00620269, 0.0058774, 0.00549714, 0.00517195, 0.00481708, 0.00447556, 0.00413199, 0.00379886, 0.0034762, 0.00314913, 0.00284838, 0.00255229, 0.00228454, 0.00201951, 0.00178857, 0.00156929, 0.00135772, 0.00117115, 0.00100593, 0.00085731, 0.00073507, 0.00062, 0.00051806, 0.00042793, 0.00036361, 0.00029552, 0.00023813, 0.00018971, 0.00015031, 0.00012512, 0.00010088, 0.00007155, 0.00005222, 0.00003577, 0.00002726, 0.00000894, 0.00000915, -0.00000638, -0.00000682, -0.0000131, -0.00001006, -0.00001248, -0.00001435, -0.00002344, -0.0000136, -0.00001425, -0.00001086, -0.0000063, 0.00000258, 0.00000169, 0.00000106, 0.00000148, -0.00000817, 0.00000243, 0.00000144, 0.00001093, 0.00001618, 0.00002091, 0.00001834, 0.0000182, 0.00002284, 0.00001619, 0.0000213, 0.000025, 0.00003535, 0.00003472, 0.00003951, 0.00004309, 0.00005558, 0.00005274, 0.00004166, 0.00005197, 0.00005572, 0.00006037, 0.00006437, 0.00006408, 0.00007078, 0.0000687, 0.00007709, 0.00007755, 0.0000874, 0.00008597, 0.0000899, 0.00009143, 0.00009581, 0.00007894, 0.00010382, 0.00009138, 0.00009977, 0.0000958, 0.00011203, 0.00012735, 0.00012372, 0.00013132, 0.0001278, 0.00013562, 0.00013278, 0.000127, 0.00013092, 0.00014118, 0.00014665, 0.00013239, 0.0001467, 0.00014544, 0.00014826, 0.00015862, 0.0001589, 0.00015052, 0.00016193, 0.00016183, 0.00017666, 0.00017634, 0.00018274, 0.0001834, 0.00018573, 0.00018249, 0.00018115, 0.00019144, 0.00018708, 0.00019388, 0.00019045, 0.00021067, 0.00019551, 0.0002043, 0.00020462, 0.00020154, 0.0002141, 0.00020841, 0.00022159, 0.00022322, 0.00022846, 0.0002242, 0.00023127, 0.00022453, 0.00022418, 0.00022716, 0.00023861, 0.00024463, 0.00023898, 0.00024646, 0.00026848, 0.0002503, 0.00026657, 0.00025943, 0.00026643, 0.00024694, 0.00026675, 0.00024928, 0.00025455, 0.00027714, 0.00026844, 0.00028045, 0.00026746, 0.00029205, 0.00028821, 0.00028009, 0.00027955, 0.00028914, 0.00028941, 0.00028376, 0.00028231, 0.00029277, 0.00030073, 0.00029925, 0.00031412, 0.0003264, 0.00032068, 0.00032787, 0.00030023, 0.00032756, 0.0003402, 0.00035477, 0.0003273, 0.00034923, 0.00034309, 0.00035172, 0.00034249, 0.00034737, 0.0003415, 0.00035132, 0.00035841, 0.00034057, 0.00035066, 0.00035675, 0.00037992, 0.00036246, 0.00038912, 0.00037307, 0.00039617, 0.00037907, 0.00037546, 0.000381, 0.00038949, 0.00041846, 0.0003972, 0.00038661, 0.00038307, 0.00039533, 0.00035532, 0.00038417, 0.00041065, 0.00038923, 0.00042529, 0.00045861, 0.00044926, 0.00042149, 0.00045738, 0.00042908, 0.0004263, 0.0004349, 0.00042999, 0.00040893, 0.00049363, 0.00045654, 0.00047574, 0.00046438, 0.00048795, 0.00047224, 0.00050204, 0.00052148, 0.00047025, 0.00047927, 0.00047166, 0.00046612, 0.00047977, 0.0005157]'

# CC_eswa
$ws.Range("H9").Value = '[-0.00045857, 0.01036226, -0.00045076, 0.00994156, -0.00423333, 0.00196317, -0.03317526, -0.01813571, -0.00042733, -0.01070245, -0.00673813, -0.00461141, 0.03161053, 0.04013503, 0.06235791, 0.0781134, 0.07659403, 0.08488783, 0.09385864, 0.09209323, 0.08889708, 0.07839895, 0.0741395, 0.06608788, 0.05534754, 0.04611076, 0.03832919, 0.03169927, 0.02609345, 0.02130389, 0.01734514, 0.01402351, 0.01130572, 0.00904298, 0.00720327, 0.00569381, 0.00449266, 0.00351752, 0.00273299, 0.0021211, 0.00164192, 0.00126519, 0.00098342, 0.00076599, 0.00059636, 0.0004688, 0.00037979, 0.00030165, 0.0002498, 0.00020469, 0.00017344, 0.00014014, 0.00011432, 0.00010334, 0.00007754, 0.00006759, 0.00005674, 0.00005001, 0.00004206, 0.00003127, 0.00003012, 0.00002981, 0.00002051, 0.0000187, 0.00000707, 0.00000691, -0.0000013, 0.00000504, -0.00000139, -0.00000335, -0.00000506, -0.00000303, -0.00001502, -0.00000249, -0.00000346, -0.00000065, -0.00000464, -0.00000713, -0.00000519, -0.00001549, -0.00001011, -0.00000665, 0.00000562, 0.00000377, 0.00001026, 0.00000783, 0.00001133, 0.00000841, 0.00001045, 0.00001635, 0.00002257, 0.00002946, 0.00003149, 0.00002839, 0.00002628, 0.00003292, 0.00003282, 0.00003351, 0.0000314, 0.00003959, 0.00005453, 0.0000511, 0.00004819, 0.00004673, 0.00004938, 0.00005855, 0.00005638, 0.00005973, 0.00006726, 0.00006463, 0.00007188, 0.00007833, 0.00007406, 0.00006386, 0.00007217, 0.00008621, 0.00007406, 0.00008446, 0.00008038, 0.00010794, 0.00008396, 0.00009449, 0.00009446, 0.0000865, 0.00009204, 0.00010327, 0.00010216, 0.00010812, 0.00010457, 0.00010496, 0.00011654, 0.00012259, 0.00012413, 0.00012004, 0.0001324, 0.00013239, 0.00012727, 0.00013725, 0.00013471, 0.00013729, 0.00014946, 0.00014538, 0.00015088, 0.00015423, 0.00015491, 0.00015813, 0.00017485, 0.00016936, 0.00017246, 0.00017287, 0.00018645, 0.00018007, 0.00018545, 0.00018232, 0.0001831, 0.00019011, 0.00018648, 0.00020316, 0.00019477, 0.00021129, 0.00020582, 0.00022908, 0.00021287, 0.0002183, 0.00023035, 0.00023181, 0.00021511, 0.00023161, 0.0002281, 0.00022974, 0.00023495, 0.00024391, 0.00025388, 0.00025916, 0.00026542, 0.00025902, 0.00027102, 0.00027029, 0.00025753, 0.00027251, 0.00028965, 0.00027765, 0.00029183, 0.00029076, 0.00028565, 0.00031312, 0.00033049, 0.00031094, 0.00030777, 0.00031899, 0.00034783, 0.00033431, 0.0003249, 0.00035073, 0.00034957, 0.00033374, 0.000342, 0.00034498, 0.00035548, 0.00036808, 0.00033374, 0.00035388, 0.00035368, 0.00036543, 0.000383, 0.00040438, 0.00040456, 0.00039418, 0.0004011, 0.00041422, 0.00041177, 0.00042423, 0.00043453, 0.00041726, 0.00042345, 0.00043817, 0.00039434, 0.00042196, 0.00041889, 0.00044659, 0.00043555, 0.00048274, 0.00048524, 0.00044422, 0.00045245, 0.00044791, 0.00047519, 0.00045082, 0.00050698, 0.00048978, 0.00045919, 0.00050543, 0.00048647, 0.00048424, 0.00051883, 0.00053224, 0.00048997, 0.00054157, 0.00055338, 0.00051653, 0.00054226, 0.00058531, 0.00048581, 0.00055085, 0.00051315, 0.00048978, 0.0005111, 0.00050743, 0.00055906, 0.00054292, 0.00060715, 0.000601, 0.00067779, 0.00059927, 0.0006055, 0.00053622]'

# CC_di
$ws.Range("H10").Value = '[-2.79166667, -3.45833333, -0.08333333, 0.75, -0.66666667, 0.91666667, 1.25, -1.625, -2.08333333, -4.79166667, -8.29166667, 0.625, 54.25, 273.70833333, 917.125, 2236, 4194.79166667, 6354.75, 8280.70833333, 9820.25, 11032.16666667, 12087.125, 13109.79166667, 14205.45833333, 15437.375, 16836.75, 18442.45833333, 20238.58333333, 22173, 24188.75, 26173.79166667, 28012.45833333, 29511.79166667, 30614.70833333, 31168.5, 31239.41666667, 30784.41666667, 29921.66666667, 28851.25, 27686.58333333, 26554.125, 25574.25, 24788.66666667, 24219.45833333, 23929.5, 23863.83333333, 24083.54166667, 24559.04166667, 25299.70833333, 26291.20833333, 27579.16666667, 29090.16666667, 30867.58333333, 32807.33333333, 34901.95833333, 37123.54166667, 39346.66666667, 41436.45833333, 43326.04166667, 44874.54166667, 46014.5, 46593.375, 46657.83333333, 46146.125, 45159.70833333, 43736.54166667, 42013.08333333, 40168.54166667, 38231.79166667, 36316.08333333, 34543.54166667, 32886.625, 31393.625, 30125.125, 28806.04166667, 28063.04166667, 27381.66666667, 26764.04166667, 26311.33333333, 26032.45833333, 25848.54166667, 25815.41666667, 25933.625, 26140.5, 26507.16666667, 26965.45833333, 27519.16666667, 28180.16666667, 28896.20833333, 29686.20833333, 30461.58333333, 31279.91666667, 32001.20833333, 32623.41666667, 33110, 33419.45833333, 33486.29166667, 33322.75, 32870.08333333, 32234.125, 31333.83333333, 30299.08333333, 29129, 27917.04166667, 26682.70833333, 25488.66666667, 24343.20833333, 23283.25, 22340.83333333, 21475.54166667, 20749.33333333, 20119.16666667, 19572.79166667, 19152.5, 18806.29166667, 18553.75, 18381.79166667, 18293.79166667, 18266.875, 18334.375, 18472.16666667, 18671.66666667, 18974.91666667, 19340.70833333, 19773.29166667, 20291.375, 20851.83333333, 21496.5, 22183.75, 22907.58333333, 23661.375, 24427.83333333, 25192.375, 25922.70833333, 26668.79166667, 27307.625, 27905.54166667, 28412.875, 28817.66666667, 29111.45833333, 29298.70833333, 29348.08333333, 29287.58333333, 29113.25, 28803.79166667, 28387.70833333, 27904.25, 27312.45833333, 26662.875, 25953.25, 25196.75, 24438.08333333, 23679.41666667, 22933.875, 22213.375, 21515.29166667, 20840.83333333, 20218.25, 19660.5, 19139.29166667, 18663.33333333, 18231.75, 17843.16666667, 17495.66666667, 17159.125, 16832.20833333, 16530.875, 16241.29166667, 15970.79166667, 15717.54166667, 15478.45833333, 15264.70833333, 15069.70833333, 14909.16666667, 14771.75, 14663.125, 14554.75, 14466.75, 14397.375, 14342, 14289.16666667, 14248.83333333, 14208.08333333, 14178.08333333, 14147, 14119.375, 14098, 14088.58333333, 14062.5, 14054.54166667, 14030.54166667, 14001.29166667, 13982.25, 13959.16666667, 13938.33333333, 13919.625, 13874.95833333, 13829.75, 13762.91666667, 13685.25, 13593.08333333, 13499.91666667, 13384.08333333, 13257.25, 13107.08333333, 12930.95833333, 12740.70833333, 12545.25, 12323.29166667, 12087.29166667, 11836.29166667, 11580.16666667, 11303.95833333, 11036.45833333, 10742.45833333, 10470.70833333, 10199.95833333, 9925, 9668.41666667, 9436.91666667, 9192.70833333, 8949.625, 8704, 8436.91666667, 8192.29166667, 7957.125, 7737.54166667, 7558.54166667, 7390.04166667, 7239.41666667, 7108.625, 6996.29166667, 6905.54166667, 6824.75, 6748.375, 6681.41666667, 6612.04166667, 6524.25, 6414.125, 6311.25, 6218, 6145.58333333, 6075.16666667, 6021.875, 5976.375, 5931.45833333, 5891.79166667, 5869.25, 5868.5, 5871.5, 5870.95833333, 5871, 5880.95833333, 5826.25, 5603.70833333, 4701.25]'

# CC_lower_wavelength_limit_for_spectra_fit
$ws.Range("H11").Value = 217
$ws.Range("H11").NumberFormat = "0.00"
$ws.Range("H11").Interior.Color = 5296274
$ws.Range("H11").HorizontalAlignment = -4131
$ws.Range("H11").VerticalAlignment = -4108

# CC_upper_wavelength_limit_for_spectra_fit - reuse H11's format to avoid creating
# a duplicate/orphan style entry, then set its own value
$ws.Range("H11").Copy()
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("H12").Value = 240
$excel.CutCopyMode = $false

# Update selection to H10 as in the saved file
$ws.Range("H10").Select()
